$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 values
$ws.Range("B74").Value = -1837
$ws.Range("C74").Value = 2547
$ws.Range("D74").Value = 4266
$ws.Range("E74").Value = 22232
$ws.Range("F74").Value = 17967
$ws.Range("G74").Value = -1719
$ws.Range("H74").Value = -4334
$ws.Range("I74").Value = -4147
$ws.Range("J74").Value = -49
$ws.Range("L74").Value = -1836
$ws.Range("M74").Value = -919
$ws.Range("N74").Value = -507
$ws.Range("O74").Value = 1445
$ws.Range("Q74").Value = -2639
$ws.Range("S74").Value = 918

# Fill in new row 75 values (and update existing E75/F75)
$ws.Range("B75").Value = -2638
$ws.Range("C75").Value = 2273
$ws.Range("D75").Value = 4174
$ws.Range("E75").Value = 23234
$ws.Range("F75").Value = 19060
$ws.Range("G75").Value = -1902
$ws.Range("H75").Value = -4922
$ws.Range("I75").Value = -4735
$ws.Range("J75").Value = 11
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = -2637
$ws.Range("M75").Value = -4795
$ws.Range("N75").Value = -3436
$ws.Range("O75").Value = -11714
$ws.Range("P75").Value = 74
$ws.Range("Q75").Value = 5835
$ws.Range("R75").Value = 4446
$ws.Range("S75").Value = -2157
$ws.Range("T75").Value = 4446
